$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column A
$ws.Range("A1").Value = "Polarity"

# New Polarity values for rows 2-68 (one per data row)
$polarity = @(-1,1,1,-1,-1,-1,-1,1,1,1,1,1,1,-1,-1,1,1,1,1,-1,1,-1,1,-1,1,1,1,-1,1,-1,-1,1,1,1,1,1,1,-1,-1,1,1,1,1,-1,1,1,-1,1,1,1,1,1,1,-1,1,-1,-1,1,-1,-1,1,1,1,1,1,-1,1)

for ($i = 0; $i -lt $polarity.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $polarity[$i]
}

# Match the updated selection recorded in the saved file
$ws.Range("H3").Select()
